$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date for rows 2-14 from 46070 to 46072
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46072
}

# Swap the data in rows 10 and 13 for columns A (Beteckning), B (Datum) and G (Area (ha))
$a10 = $ws.Cells.Item(10, 1).Value2
$b10 = $ws.Cells.Item(10, 2).Value2
$g10 = $ws.Cells.Item(10, 7).Value2

$a13 = $ws.Cells.Item(13, 1).Value2
$b13 = $ws.Cells.Item(13, 2).Value2
$g13 = $ws.Cells.Item(13, 7).Value2

$ws.Cells.Item(10, 1).Value2 = $a13
$ws.Cells.Item(10, 2).Value2 = $b13
$ws.Cells.Item(10, 7).Value2 = $g13

$ws.Cells.Item(13, 1).Value2 = $a10
$ws.Cells.Item(13, 2).Value2 = $b10
$ws.Cells.Item(13, 7).Value2 = $g10
